$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit moves the row for "Fattoria San Leolino" from position 2 down to
# position 14 (everything between shifts up by one row); rows 15-26 are untouched.
# We rewrite the affected block (rows 2-14) with the final name/url content and
# then rebuild the worksheet hyperlinks (the COM engine only exposes a whole-sheet
# Hyperlinks.Delete(), so every hyperlink is re-added, including the ones on rows
# that did not otherwise change).

$ws.Range("A2").Value = "Fattoria il Casalone"
$ws.Range("B2").Value = "http://www.agriturismocasalone.com"

$ws.Range("A3").Value = "Azienda Agricola Barbialla Srl Società Agricola - Castello di Barbialla"
$ws.Range("B3").Value = "https://www.aziendaagricolabarbialla.it"

$ws.Range("A4").Value = "Elleci Riviera"
$ws.Range("B4").Value = "https://elleciriviera.it/"

$ws.Range("A5").Value = "Azienda Agricola Ulivelli"
$ws.Range("B5").Value = "http://www.agricolaulivelli.it"

$ws.Range("A6").Value = "Fattoria Ruschi Noceti"
$ws.Range("B6").Value = "http://www.fattoriaruschinoceti.com"

$ws.Range("A7").Value = "Saniscope-Chimica"
$ws.Range("B7").Value = "http://www.saniscope-chimica.it/"

$ws.Range("A8").Value = "Frantoio Carmignani Umberto"
$ws.Range("B8").Value = "N/A"
$ws.Range("B8").Style = "Normal"

$ws.Range("A9").Value = "Castello di Gabbiano"
$ws.Range("B9").Value = "http://www.castellogabbiano.it"

$ws.Range("A10").Value = "Tenute Silvio Nardi"
$ws.Range("B10").Value = "https://www.tenutenardi.com/"

$ws.Range("A11").Value = "Rocca di Castagnoli Società Agricola"
$ws.Range("B11").Value = "https://www.roccadicastagnoli.com/"

$ws.Range("A12").Value = "Archa"
$ws.Range("B12").Value = "https://www.archa.it/"

$ws.Range("A13").Value = "Ortofrutta Branchi"
$ws.Range("B13").Value = "N/A"
$ws.Range("B13").Style = "Normal"

$ws.Range("A14").Value = "Fattoria San Leolino"
$ws.Range("B14").Value = "http://www.fattoriasanleolino.com"

# Rebuild hyperlinks for every row that has a real URL (i.e. not "N/A").
$ws.Range("A1").Hyperlinks.Delete()

$h = $ws.Hyperlinks.Add($ws.Range("B2"), "http://www.agriturismocasalone.com")
$h.Address = "http://www.agriturismocasalone.com"
$ws.Range("B2").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B3"), "https://www.aziendaagricolabarbialla.it")
$h.Address = "https://www.aziendaagricolabarbialla.it"
$ws.Range("B3").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B4"), "https://elleciriviera.it/")
$h.Address = "https://elleciriviera.it/"
$ws.Range("B4").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B5"), "http://www.agricolaulivelli.it")
$h.Address = "http://www.agricolaulivelli.it"
$ws.Range("B5").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B6"), "http://www.fattoriaruschinoceti.com")
$h.Address = "http://www.fattoriaruschinoceti.com"
$ws.Range("B6").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B7"), "http://www.saniscope-chimica.it/")
$h.Address = "http://www.saniscope-chimica.it/"
$ws.Range("B7").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B9"), "http://www.castellogabbiano.it")
$h.Address = "http://www.castellogabbiano.it"
$ws.Range("B9").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B10"), "https://www.tenutenardi.com/")
$h.Address = "https://www.tenutenardi.com/"
$ws.Range("B10").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B11"), "https://www.roccadicastagnoli.com/")
$h.Address = "https://www.roccadicastagnoli.com/"
$ws.Range("B11").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B12"), "https://www.archa.it/")
$h.Address = "https://www.archa.it/"
$ws.Range("B12").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B14"), "http://www.fattoriasanleolino.com")
$h.Address = "http://www.fattoriasanleolino.com"
$ws.Range("B14").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B15"), "https://www.aziendagricolamagnanigigliola.com")
$h.Address = "https://www.aziendagricolamagnanigigliola.com"
$ws.Range("B15").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B17"), "http://www.tenutalafortuna.it")
$h.Address = "http://www.tenutalafortuna.it"
$ws.Range("B17").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B19"), "https://www.impresadipulizietotalclean.it")
$h.Address = "https://www.impresadipulizietotalclean.it"
$ws.Range("B19").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B20"), "https://www.lenziagricoltura.it")
$h.Address = "https://www.lenziagricoltura.it"
$ws.Range("B20").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B22"), "http://www.terredifirenze.it")
$h.Address = "http://www.terredifirenze.it"
$ws.Range("B22").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B23"), "http://www.impresadipulizielafonte.it")
$h.Address = "http://www.impresadipulizielafonte.it"
$ws.Range("B23").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B24"), "https://www.aziendaagricolacanciulle.it")
$h.Address = "https://www.aziendaagricolacanciulle.it"
$ws.Range("B24").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B25"), "https://www.impresapulizietasselli.it")
$h.Address = "https://www.impresapulizietasselli.it"
$ws.Range("B25").Style = "Hyperlink"

$h = $ws.Hyperlinks.Add($ws.Range("B26"), "https://www.blitzservice.it")
$h.Address = "https://www.blitzservice.it"
$ws.Range("B26").Style = "Hyperlink"

